$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = 45967
$ws.Range("B23").Value = 604
$ws.Range("C23").Value = 14
$ws.Range("D23").Value = 590

$ws.Range("A23:D23").Select()
